# cryptos.xlsx periodic price-refresh update (GitHub Actions style commit).
#
# The sheet stores every Price/Volume(1h) cell as literal text (inlineStr in
# the OOXML), never as a real number, so prices like "1.00" or "0.998" must
# round-trip as text, not get auto-coerced into numeric cells by Excel's
# smart-entry heuristics. Cells whose new text cannot be parsed as a plain
# number (e.g. thousands-grouped "68.054.22", names, URLs, the padded
# "  -2.22%  " percentages) are safe to assign directly. Cells whose new
# text WOULD parse as a plain number (e.g. "0.998", "7.27") are routed
# through a brief Text-format -> assign -> ClearFormats dance so the final
# cell keeps its original (default) style while the stored value stays text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Auto-generated body: direct text assignments for non-ambiguous strings
$ws.Range('D2').Value = '68.054.22'
$ws.Range('E2').Value = '  -2.22%  '
$ws.Range('D3').Value = '3.790.02'
$ws.Range('E3').Value = '  +2.65%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('E5').Value = '  -3.58%  '
$ws.Range('E6').Value = '  -4.16%  '
$ws.Range('D7').Value = '3.787.05'
$ws.Range('E7').Value = '  +2.57%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('E10').Value = '  -2.59%  '
$ws.Range('E11').Value = '  +0.81%  '
$ws.Range('E12').Value = '  -2.29%  '
$ws.Range('E13').Value = '  -3.78%  '
$ws.Range('E14').Value = '  -3.22%  '
$ws.Range('D15').Value = '4.417.83'
$ws.Range('E15').Value = '  +2.73%  '
$ws.Range('D16').Value = '3.784.30'
$ws.Range('E16').Value = '  +2.83%  '
$ws.Range('D17').Value = '68.225.67'
$ws.Range('E17').Value = '  -1.89%  '
$ws.Range('B18').Value = 'Polkadot'
$ws.Range('C18').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('E18').Value = '  -3.49%  '
$ws.Range('B19').Value = 'TRON'
$ws.Range('C19').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('E19').Value = '  -4.43%  '
$ws.Range('E20').Value = '  -1.59%  '
$ws.Range('E21').Value = '  -2.00%  '
$ws.Range('E22').Value = '  +3.71%  '
$ws.Range('E23').Value = '  +1.90%  '
$ws.Range('E24').Value = '  -0.16%  '
$ws.Range('E25').Value = '  -5.19%  '
$ws.Range('E26').Value = '  +8.13%  '
$ws.Range('E27').Value = '  -5.09%  '
$ws.Range('E28').Value = '  -8.50%  '
$ws.Range('E29').Value = '  +0.08%  '
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('E31').Value = '  +0.37%  '
$ws.Range('E32').Value = '  +7.41%  '
$ws.Range('E33').Value = '  -3.81%  '
$ws.Range('E34').Value = '  -2.46%  '
$ws.Range('E35').Value = '  +0.05%  '
$ws.Range('E36').Value = '  -2.26%  '
$ws.Range('E37').Value = '  -2.88%  '
$ws.Range('E38').Value = '  -1.68%  '
$ws.Range('E39').Value = '  -4.06%  '
$ws.Range('E40').Value = '  +5.37%  '
$ws.Range('E41').Value = '  -1.62%  '
$ws.Range('E42').Value = '  -1.90%  '
$ws.Range('E43').Value = '  -1.66%  '
$ws.Range('E44').Value = '  -2.33%  '
$ws.Range('E45').Value = '  -4.52%  '
$ws.Range('D46').Value = '2.846.44'
$ws.Range('E46').Value = '  -3.17%  '
$ws.Range('B47').Value = 'USDe'
$ws.Range('C47').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('E47').Value = '  +0.06%  '
$ws.Range('B48').Value = 'VeChain'
$ws.Range('C48').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('E48').Value = '  -1.99%  '
$ws.Range('E49').Value = '  +0.81%  '
$ws.Range('E50').Value = '  -2.44%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('E51').Value = '  +10.33%  '

# Numeric-looking strings must be forced to store as text (matches source inlineStr)
# by briefly tagging the cell as Text-formatted, then clearing the format back
# to General so no stray style survives on the cell.
$forceTextCells = @(
    @{Ref='D5'; Val='593.36'}
    @{Ref='D6'; Val='171.78'}
    @{Ref='D9'; Val='0.529'}
    @{Ref='D10'; Val='0.160'}
    @{Ref='D11'; Val='6.30'}
    @{Ref='D13'; Val='38.34'}
    @{Ref='D14'; Val='0.0000245'}
    @{Ref='D18'; Val='7.27'}
    @{Ref='D19'; Val='0.116'}
    @{Ref='D20'; Val='16.05'}
    @{Ref='D21'; Val='488.63'}
    @{Ref='D22'; Val='9.46'}
    @{Ref='D23'; Val='0.732'}
    @{Ref='D24'; Val='85.96'}
    @{Ref='D25'; Val='2.37'}
    @{Ref='D26'; Val='0.0000139'}
    @{Ref='D28'; Val='10.22'}
    @{Ref='D30'; Val='2.93'}
    @{Ref='D31'; Val='2.44'}
    @{Ref='D32'; Val='32.35'}
    @{Ref='D33'; Val='7.63'}
    @{Ref='D35'; Val='0.998'}
    @{Ref='D36'; Val='1.02'}
    @{Ref='D37'; Val='5.86'}
    @{Ref='D39'; Val='0.326'}
    @{Ref='D40'; Val='451.39'}
    @{Ref='D41'; Val='49.16'}
    @{Ref='D43'; Val='2.87'}
    @{Ref='D44'; Val='8.36'}
    @{Ref='D45'; Val='41.71'}
    @{Ref='D47'; Val='1.00'}
    @{Ref='D48'; Val='0.0352'}
    @{Ref='D49'; Val='137.92'}
    @{Ref='D50'; Val='26.70'}
    @{Ref='D51'; Val='23.52'}
)

foreach ($item in $forceTextCells) {
    $c = $ws.Range($item.Ref)
    $c.NumberFormat = "@"
    $c.Value = $item.Val
    $c.ClearFormats()
}
